$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
$ws.Range("D2").Value = 44827
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 30000
$ws.Range("L2").Value = 31000
$ws.Range("M2").Value = 30500
$ws.Range("P2").Value = 1220

# Row 4 updates
$ws.Range("D4").Value = 44379
$ws.Range("J4").Value = 240
$ws.Range("K4").Value = 31000
$ws.Range("L4").Value = 32000
$ws.Range("M4").Value = 31500
$ws.Range("P4").Value = 1260
